$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 21.750601
$ws.Cells.Item(2, 8).Value = 65.251803
$ws.Cells.Item(2, 9).Value = 0.3612485837751334
$ws.Cells.Item(2, 10).Value = 0.3773020033645041
$ws.Cells.Item(2, 13).Value = 14.23612
$ws.Cells.Item(2, 14).Value = 42.70836
$ws.Cells.Item(2, 15).Value = 0.07600219901258977
$ws.Cells.Item(2, 16).Value = 0.09845490883293605
$ws.Cells.Item(2, 17).Value = 309.64416590812
$ws.Cells.Item(2, 18).Value = 2786.79749317308
$ws.Cells.Item(2, 19).Value = 0.0274556867570939
$ws.Cells.Item(2, 20).Value = 0.03714723434373638
$ws.Cells.Item(3, 7).Value = 21.750601
$ws.Cells.Item(3, 8).Value = 65.251803
$ws.Cells.Item(3, 9).Value = 0.3612485837751334
$ws.Cells.Item(3, 10).Value = 0.3773020033645041
$ws.Cells.Item(3, 15).Value = 0.1644984167819043
$ws.Cells.Item(3, 16).Value = 0.2130948424892534
$ws.Cells.Item(3, 17).Value = 670.190806573919
$ws.Cells.Item(3, 18).Value = 6031.717259165271
$ws.Cells.Item(3, 19).Value = 0.05942482009571456
$ws.Cells.Item(3, 20).Value = 0.08040111097783874
$ws.Cells.Item(4, 7).Value = 21.750601
$ws.Cells.Item(4, 8).Value = 65.251803
$ws.Cells.Item(4, 9).Value = 0.3612485837751334
$ws.Cells.Item(4, 10).Value = 0.3773020033645041
$ws.Cells.Item(4, 13).Value = 5.407681
$ws.Cells.Item(4, 14).Value = 16.223043
$ws.Cells.Item(4, 15).Value = 0.02886992014387351
$ws.Cells.Item(4, 16).Value = 0.03739872520410058
$ws.Cells.Item(4, 17).Value = 117.620311766281
$ws.Cells.Item(4, 18).Value = 1058.582805896529
$ws.Cells.Item(4, 19).Value = 0.0104292177656755
$ws.Cells.Item(4, 20).Value = 0.01411061394278572
$ws.Cells.Item(5, 7).Value = 21.750601
$ws.Cells.Item(5, 8).Value = 65.251803
$ws.Cells.Item(5, 9).Value = 0.3612485837751334
$ws.Cells.Item(5, 10).Value = 0.3773020033645041
$ws.Cells.Item(5, 13).Value = 128.149857
$ws.Cells.Item(5, 14).Value = 256.299714
$ws.Cells.Item(5, 15).Value = 0.6841520677789258
$ws.Cells.Item(5, 16).Value = 0.5908436890523912
$ws.Cells.Item(5, 17).Value = 2787.336407814057
$ws.Cells.Item(5, 18).Value = 16724.01844688434
$ws.Cells.Item(5, 19).Value = 0.247148965571966
$ws.Cells.Item(5, 20).Value = 0.2229265075547413
$ws.Cells.Item(6, 7).Value = 21.750601
$ws.Cells.Item(6, 8).Value = 65.251803
$ws.Cells.Item(6, 9).Value = 0.3612485837751334
$ws.Cells.Item(6, 10).Value = 0.3773020033645041
$ws.Cells.Item(6, 13).Value = 8.705771666666665
$ws.Cells.Item(6, 14).Value = 26.117315
$ws.Cells.Item(6, 15).Value = 0.04647739628270661
$ws.Cells.Item(6, 16).Value = 0.06020783442131875
$ws.Cells.Item(6, 17).Value = 189.3557659187716
$ws.Cells.Item(6, 18).Value = 1704.201893268945
$ws.Cells.Item(6, 19).Value = 0.01678989358468341
$ws.Cells.Item(6, 20).Value = 0.02271653654540191
$ws.Cells.Item(7, 9).Value = 0.2797740820980411
$ws.Cells.Item(7, 10).Value = 0.2922068802649305
$ws.Cells.Item(7, 13).Value = 14.23612
$ws.Cells.Item(7, 14).Value = 42.70836
$ws.Cells.Item(7, 15).Value = 0.07600219901258977
$ws.Cells.Item(7, 16).Value = 0.09845490883293605
$ws.Cells.Item(7, 17).Value = 239.80830980332
$ws.Cells.Item(7, 18).Value = 2158.27478822988
$ws.Cells.Item(7, 19).Value = 0.02126344546617994
$ws.Cells.Item(7, 20).Value = 0.02876920175684039
$ws.Cells.Item(8, 9).Value = 0.2797740820980411
$ws.Cells.Item(8, 10).Value = 0.2922068802649305
$ws.Cells.Item(8, 15).Value = 0.1644984167819043
$ws.Cells.Item(8, 16).Value = 0.2130948424892534
$ws.Cells.Item(8, 19).Value = 0.04602239356173827
$ws.Cells.Item(8, 20).Value = 0.06226777912433148
$ws.Cells.Item(9, 9).Value = 0.2797740820980411
$ws.Cells.Item(9, 10).Value = 0.2922068802649305
$ws.Cells.Item(9, 13).Value = 5.407681
$ws.Cells.Item(9, 14).Value = 16.223043
$ws.Cells.Item(9, 15).Value = 0.02886992014387351
$ws.Cells.Item(9, 16).Value = 0.03739872520410058
$ws.Cells.Item(9, 17).Value = 91.092716313541
$ws.Cells.Item(9, 18).Value = 819.8344468218691
$ws.Cells.Item(9, 19).Value = 0.008077055408495955
$ws.Cells.Item(9, 20).Value = 0.01092816481777566
$ws.Cells.Item(10, 9).Value = 0.2797740820980411
$ws.Cells.Item(10, 10).Value = 0.2922068802649305
$ws.Cells.Item(10, 13).Value = 128.149857
$ws.Cells.Item(10, 14).Value = 256.299714
$ws.Cells.Item(10, 15).Value = 0.6841520677789258
$ws.Cells.Item(10, 16).Value = 0.5908436890523912
$ws.Cells.Item(10, 17).Value = 2158.692158306277
$ws.Cells.Item(10, 18).Value = 12952.15294983766
$ws.Cells.Item(10, 19).Value = 0.1914080167783257
$ws.Cells.Item(10, 20).Value = 0.1726485911022219
$ws.Cells.Item(11, 9).Value = 0.2797740820980411
$ws.Cells.Item(11, 10).Value = 0.2922068802649305
$ws.Cells.Item(11, 13).Value = 8.705771666666665
$ws.Cells.Item(11, 14).Value = 26.117315
$ws.Cells.Item(11, 15).Value = 0.04647739628270661
$ws.Cells.Item(11, 16).Value = 0.06020783442131875
$ws.Cells.Item(11, 17).Value = 146.6492547770717
$ws.Cells.Item(11, 18).Value = 1319.843292993645
$ws.Cells.Item(11, 19).Value = 0.01300317088330115
$ws.Cells.Item(11, 20).Value = 0.01759314346376105
$ws.Cells.Item(12, 7).Value = 8.938416999999999
$ws.Cells.Item(12, 8).Value = 26.815251
$ws.Cells.Item(12, 9).Value = 0.1484552303838214
$ws.Cells.Item(12, 10).Value = 0.1550523887136425
$ws.Cells.Item(12, 13).Value = 14.23612
$ws.Cells.Item(12, 14).Value = 42.70836
$ws.Cells.Item(12, 15).Value = 0.07600219901258977
$ws.Cells.Item(12, 16).Value = 0.09845490883293605
$ws.Cells.Item(12, 17).Value = 127.24837702204
$ws.Cells.Item(12, 18).Value = 1145.23539319836
$ws.Cells.Item(12, 19).Value = 0.01128292396409106
$ws.Cells.Item(12, 20).Value = 0.01526566879513063
$ws.Cells.Item(13, 7).Value = 8.938416999999999
$ws.Cells.Item(13, 8).Value = 26.815251
$ws.Cells.Item(13, 9).Value = 0.1484552303838214
$ws.Cells.Item(13, 10).Value = 0.1550523887136425
$ws.Cells.Item(13, 15).Value = 0.1644984167819043
$ws.Cells.Item(13, 16).Value = 0.2130948424892534
$ws.Cells.Item(13, 17).Value = 275.4151436424229
$ws.Cells.Item(13, 18).Value = 2478.736292781807
$ws.Cells.Item(13, 19).Value = 0.02442065036113148
$ws.Cells.Item(13, 20).Value = 0.03304086435051613
$ws.Cells.Item(14, 7).Value = 8.938416999999999
$ws.Cells.Item(14, 8).Value = 26.815251
$ws.Cells.Item(14, 9).Value = 0.1484552303838214
$ws.Cells.Item(14, 10).Value = 0.1550523887136425
$ws.Cells.Item(14, 13).Value = 5.407681
$ws.Cells.Item(14, 14).Value = 16.223043
$ws.Cells.Item(14, 15).Value = 0.02886992014387351
$ws.Cells.Item(14, 16).Value = 0.03739872520410058
$ws.Cells.Item(14, 17).Value = 48.336107780977
$ws.Cells.Item(14, 18).Value = 435.024970028793
$ws.Cells.Item(14, 19).Value = 0.004285890646121268
$ws.Cells.Item(14, 20).Value = 0.005798761677740901
$ws.Cells.Item(15, 7).Value = 8.938416999999999
$ws.Cells.Item(15, 8).Value = 26.815251
$ws.Cells.Item(15, 9).Value = 0.1484552303838214
$ws.Cells.Item(15, 10).Value = 0.1550523887136425
$ws.Cells.Item(15, 13).Value = 128.149857
$ws.Cells.Item(15, 14).Value = 256.299714
$ws.Cells.Item(15, 15).Value = 0.6841520677789258
$ws.Cells.Item(15, 16).Value = 0.5908436890523912
$ws.Cells.Item(15, 17).Value = 1145.456860356369
$ws.Cells.Item(15, 18).Value = 6872.741162138213
$ws.Cells.Item(15, 19).Value = 0.1015659528396882
$ws.Cells.Item(15, 20).Value = 0.09161172534395386
$ws.Cells.Item(16, 7).Value = 8.938416999999999
$ws.Cells.Item(16, 8).Value = 26.815251
$ws.Cells.Item(16, 9).Value = 0.1484552303838214
$ws.Cells.Item(16, 10).Value = 0.1550523887136425
$ws.Cells.Item(16, 13).Value = 8.705771666666665
$ws.Cells.Item(16, 14).Value = 26.117315
$ws.Cells.Item(16, 15).Value = 0.04647739628270661
$ws.Cells.Item(16, 16).Value = 0.06020783442131875
$ws.Cells.Item(16, 17).Value = 77.81581746345165
$ws.Cells.Item(16, 18).Value = 700.3423571710649
$ws.Cells.Item(16, 19).Value = 0.006899812572789375
$ws.Cells.Item(16, 20).Value = 0.009335368546300937
$ws.Cells.Item(17, 7).Value = 7.6853705
$ws.Cells.Item(17, 8).Value = 15.370741
$ws.Cells.Item(17, 9).Value = 0.1276437928732263
$ws.Cells.Item(17, 10).Value = 0.08887741190073968
$ws.Cells.Item(17, 13).Value = 14.23612
$ws.Cells.Item(17, 14).Value = 42.70836
$ws.Cells.Item(17, 15).Value = 0.07600219901258977
$ws.Cells.Item(17, 16).Value = 0.09845490883293605
$ws.Cells.Item(17, 17).Value = 109.40985668246
$ws.Cells.Item(17, 18).Value = 656.45914009476
$ws.Cells.Item(17, 19).Value = 0.009701208948672735
$ws.Cells.Item(17, 20).Value = 0.00875041748599463
$ws.Cells.Item(18, 7).Value = 7.6853705
$ws.Cells.Item(18, 8).Value = 15.370741
$ws.Cells.Item(18, 9).Value = 0.1276437928732263
$ws.Cells.Item(18, 10).Value = 0.08887741190073968
$ws.Cells.Item(18, 15).Value = 0.1644984167819043
$ws.Cells.Item(18, 16).Value = 0.2130948424892534
$ws.Cells.Item(18, 17).Value = 236.8056245532895
$ws.Cells.Item(18, 18).Value = 1420.833747319737
$ws.Cells.Item(18, 19).Value = 0.02099720183968305
$ws.Cells.Item(18, 20).Value = 0.01893931808984061
$ws.Cells.Item(19, 7).Value = 7.6853705
$ws.Cells.Item(19, 8).Value = 15.370741
$ws.Cells.Item(19, 9).Value = 0.1276437928732263
$ws.Cells.Item(19, 10).Value = 0.08887741190073968
$ws.Cells.Item(19, 13).Value = 5.407681
$ws.Cells.Item(19, 14).Value = 16.223043
$ws.Cells.Item(19, 15).Value = 0.02886992014387351
$ws.Cells.Item(19, 16).Value = 0.03739872520410058
$ws.Cells.Item(19, 17).Value = 41.5600320308105
$ws.Cells.Item(19, 18).Value = 249.360192184863
$ws.Cells.Item(19, 19).Value = 0.003685066107111174
$ws.Cells.Item(19, 20).Value = 0.003323901904527422
$ws.Cells.Item(20, 7).Value = 7.6853705
$ws.Cells.Item(20, 8).Value = 15.370741
$ws.Cells.Item(20, 9).Value = 0.1276437928732263
$ws.Cells.Item(20, 10).Value = 0.08887741190073968
$ws.Cells.Item(20, 13).Value = 128.149857
$ws.Cells.Item(20, 14).Value = 256.299714
$ws.Cells.Item(20, 15).Value = 0.6841520677789258
$ws.Cells.Item(20, 16).Value = 0.5908436890523912
$ws.Cells.Item(20, 17).Value = 984.8791305670185
$ws.Cells.Item(20, 18).Value = 3939.516522268074
$ws.Cells.Item(20, 19).Value = 0.0873277648333627
$ws.Cells.Item(20, 20).Value = 0.05251265792086193
$ws.Cells.Item(21, 7).Value = 7.6853705
$ws.Cells.Item(21, 8).Value = 15.370741
$ws.Cells.Item(21, 9).Value = 0.1276437928732263
$ws.Cells.Item(21, 10).Value = 0.08887741190073968
$ws.Cells.Item(21, 13).Value = 8.705771666666665
$ws.Cells.Item(21, 14).Value = 26.117315
$ws.Cells.Item(21, 15).Value = 0.04647739628270661
$ws.Cells.Item(21, 16).Value = 0.06020783442131875
$ws.Cells.Item(21, 17).Value = 66.90708074673583
$ws.Cells.Item(21, 18).Value = 401.442484480415
$ws.Cells.Item(21, 19).Value = 0.005932551144396662
$ws.Cells.Item(21, 20).Value = 0.005351116499515079
$ws.Cells.Item(22, 7).Value = 4.990062666666667
$ws.Cells.Item(22, 8).Value = 14.970188
$ws.Cells.Item(22, 9).Value = 0.08287831086977776
$ws.Cells.Item(22, 10).Value = 0.08656131575618316
$ws.Cells.Item(22, 13).Value = 14.23612
$ws.Cells.Item(22, 14).Value = 42.70836
$ws.Cells.Item(22, 15).Value = 0.07600219901258977
$ws.Cells.Item(22, 16).Value = 0.09845490883293605
$ws.Cells.Item(22, 17).Value = 71.03913093018667
$ws.Cells.Item(22, 18).Value = 639.35217837168
$ws.Cells.Item(22, 19).Value = 0.006298933876552131
$ws.Cells.Item(22, 20).Value = 0.008522386451234003
$ws.Cells.Item(23, 7).Value = 4.990062666666667
$ws.Cells.Item(23, 8).Value = 14.970188
$ws.Cells.Item(23, 9).Value = 0.08287831086977776
$ws.Cells.Item(23, 10).Value = 0.08656131575618316
$ws.Cells.Item(23, 15).Value = 0.1644984167819043
$ws.Cells.Item(23, 16).Value = 0.2130948424892534
$ws.Cells.Item(23, 17).Value = 153.7564007278573
$ws.Cells.Item(23, 18).Value = 1383.807606550716
$ws.Cells.Item(23, 19).Value = 0.01363335092363693
$ws.Cells.Item(23, 20).Value = 0.01844576994672638
$ws.Cells.Item(24, 7).Value = 4.990062666666667
$ws.Cells.Item(24, 8).Value = 14.970188
$ws.Cells.Item(24, 9).Value = 0.08287831086977776
$ws.Cells.Item(24, 10).Value = 0.08656131575618316
$ws.Cells.Item(24, 13).Value = 5.407681
$ws.Cells.Item(24, 14).Value = 16.223043
$ws.Cells.Item(24, 15).Value = 0.02886992014387351
$ws.Cells.Item(24, 16).Value = 0.03739872520410058
$ws.Cells.Item(24, 17).Value = 26.98466707134267
$ws.Cells.Item(24, 18).Value = 242.862003642084
$ws.Cells.Item(24, 19).Value = 0.002392690216469608
$ws.Cells.Item(24, 20).Value = 0.003237282861270876
$ws.Cells.Item(25, 7).Value = 4.990062666666667
$ws.Cells.Item(25, 8).Value = 14.970188
$ws.Cells.Item(25, 9).Value = 0.08287831086977776
$ws.Cells.Item(25, 10).Value = 0.08656131575618316
$ws.Cells.Item(25, 13).Value = 128.149857
$ws.Cells.Item(25, 14).Value = 256.299714
$ws.Cells.Item(25, 15).Value = 0.6841520677789258
$ws.Cells.Item(25, 16).Value = 0.5908436890523912
$ws.Cells.Item(25, 17).Value = 639.475817154372
$ws.Cells.Item(25, 18).Value = 3836.854902926232
$ws.Cells.Item(25, 19).Value = 0.05670136775558308
$ws.Cells.Item(25, 20).Value = 0.05114420713061213
$ws.Cells.Item(26, 7).Value = 4.990062666666667
$ws.Cells.Item(26, 8).Value = 14.970188
$ws.Cells.Item(26, 9).Value = 0.08287831086977776
$ws.Cells.Item(26, 10).Value = 0.08656131575618316
$ws.Cells.Item(26, 13).Value = 8.705771666666665
$ws.Cells.Item(26, 14).Value = 26.117315
$ws.Cells.Item(26, 15).Value = 0.04647739628270661
$ws.Cells.Item(26, 16).Value = 0.06020783442131875
$ws.Cells.Item(26, 17).Value = 43.44234617835777
$ws.Cells.Item(26, 18).Value = 390.98111560522
$ws.Cells.Item(26, 19).Value = 0.003851968097536012
$ws.Cells.Item(26, 20).Value = 0.005211669366339765
